# Horarios actualizados Linea 141 - 947
# Updates the three schedule sheets (LP1912, LP1912-215, 6203-6173) with the
# latest scrape timestamp (04:53:24) and refreshed arrival rows.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# Sheet 1: LP1912
# ---------------------------------------------------------------------------
$ws1 = $wb.Worksheets.Item("LP1912")

$ws1.Range("A2").Value = "Última actualización: 04:53:24"
$ws1.Range("A3").Value = "Total filas: 12"

$sheet1Data = @(
    @("04:53:24", "04:53", "11_ETCHEVERRY", 0, "LP1912"),
    @("04:53:24", "05:17", "17_ROMERO", 24, "LP1912"),
    @("04:53:24", "05:22", "23_HERNANDEZ", 29, "LP1912"),
    @("04:53:24", "05:44", "14_ABASTO", 51, "LP1912"),
    @("04:53:24", "05:47", "17_ROMERO", 54, "LP1912"),
    @("04:53:24", "06:01", "16_SANTA ANA", 68, "LP1912"),
    @("04:53:24", "06:09", "10_OLMOS", 76, "LP1912"),
    @("04:53:24", "06:16", "215A_EL PATO", 83, "LP1912"),
    @("04:53:24", "06:30", "23_HERNANDEZ", 97, "LP1912"),
    @("04:53:24", "06:34", "11_ETCHEVERRY", 101, "LP1912"),
    @("04:53:24", "06:39", "17X38_ROMERO", 106, "LP1912"),
    @("04:53:24", "06:41", "16_SANTA ANA", 108, "LP1912")
)

$row = 6
foreach ($r in $sheet1Data) {
    $ws1.Cells($row, 1).Value = $r[0]
    $ws1.Cells($row, 2).Value = $r[1]
    $ws1.Cells($row, 3).Value = $r[2]
    $ws1.Cells($row, 4).Value = $r[3]
    $ws1.Cells($row, 5).Value = $r[4]
    $row = $row + 1
}

# Old sheet had rows up to 20; new data only goes through row 17, so clear
# the trailing rows that no longer exist.
$ws1.Range("A18:E20").ClearContents()

# ---------------------------------------------------------------------------
# Sheet 2: LP1912-215
# ---------------------------------------------------------------------------
$ws2 = $wb.Worksheets.Item("LP1912-215")

$ws2.Range("A2").Value = "Última actualización: 04:53:24"
$ws2.Range("A3").Value = "Total filas: 1"

$ws2.Cells(6, 1).Value = "04:53:24"
$ws2.Cells(6, 2).Value = "06:16"
$ws2.Cells(6, 3).Value = "215A_EL PATO"
$ws2.Cells(6, 4).Value = 83
$ws2.Cells(6, 5).Value = "LP1912"

# Old sheet had a row 7; new data only has a single row, so clear it.
$ws2.Range("A7:E7").ClearContents()

# ---------------------------------------------------------------------------
# Sheet 3: 6203-6173
# ---------------------------------------------------------------------------
$ws3 = $wb.Worksheets.Item("6203-6173")

$ws3.Range("A2").Value = "Última actualización: 04:53:24"

Write-Output "Schedules refreshed."
